$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("B3").Value = 115.184878
$ws.Range("D3").Value = 1.278897
$ws.Range("E3").Value = 0.280378

# Row 4
$ws.Range("B4").Value = 10042.337357
$ws.Range("C4").Value = 223

# Row 5
$ws.Range("G5").Value = 1.944444
$ws.Range("H5").Value = -0.986341
$ws.Range("I5").Value = 4.87523
$ws.Range("J5").Value = 0.262836

# Row 6
$ws.Range("G6").Value = 1.140174
$ws.Range("H6").Value = -2.015661
$ws.Range("I6").Value = 4.29601
$ws.Range("J6").Value = 0.670782

# Row 7
$ws.Range("G7").Value = -0.80427
$ws.Range("H7").Value = -3.190115
$ws.Range("I7").Value = 1.581575
$ws.Range("J7").Value = 0.706256
